# Insert a new translation row "errors.NO_MODELS" just above the existing
# "button.generate" row (row 415), pushing every row from 415..451 down by
# one (to 416..452). This mirrors the diff, which inserts a brand-new key
# at A415/B415 (no Chinese translation) and leaves all subsequent rows'
# content untouched other than the row-number shift.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 415:451 down to 416:452, leaving row 415 blank for the new entry.
$ws.Rows.Item(415).Insert()

# Populate the newly inserted row with the new translation key.
$ws.Cells.Item(415, 1).Value2 = "errors.NO_MODELS"
$ws.Cells.Item(415, 2).Value2 = "No models available. Please add a model by clicking the chart in the main view."
